# Apply scheduled-runner market-data refresh to the profit-calc sheets.
# Each sheet has columns H..N holding cached market prices / profit figures
# (plain cached values, no formulas) that get refreshed by the runner.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 1404343.2
$ws.Range("J51").Value = 1491698
$ws.Range("L51").Value = 1491698
$ws.Range("N51").Value = -1492666
$ws.Range("H86").Value = 11827550
$ws.Range("J86").Value = 20103404
$ws.Range("L86").Value = 20103404
$ws.Range("N86").Value = -20105650
$ws.Range("H89").Value = 11827550
$ws.Range("J89").Value = 20103404
$ws.Range("L89").Value = 100517020
$ws.Range("N89").Value = -100528252
$ws.Range("H100").Value = 7937.375
$ws.Range("I100").Value = 11598.6
$ws.Range("J100").Value = 1835.3334
$ws.Range("K100").Value = 11598.6
$ws.Range("L100").Value = 1835.3334
$ws.Range("M100").Value = -11057.6
$ws.Range("N100").Value = -2917.3334
$ws.Range("H132").Value = 3505.7458
$ws.Range("I132").Value = 3302.9412
$ws.Range("J132").Value = 4798.625
$ws.Range("K132").Value = 9908.8236
$ws.Range("L132").Value = 14395.875
$ws.Range("M132").Value = -7378.8236
$ws.Range("N132").Value = -19455.875
$ws.Range("H137").Value = 38464296
$ws.Range("I137").Value = 200001980
$ws.Range("J137").Value = 2942.8096
$ws.Range("K137").Value = 600005940
$ws.Range("L137").Value = 8828.4288
$ws.Range("M137").Value = -600003390
$ws.Range("N137").Value = -13928.4288
$ws.Range("H138").Value = 1883.3673
$ws.Range("I138").Value = 1399.0344
$ws.Range("J138").Value = 2585.65
$ws.Range("K138").Value = 4197.1032
$ws.Range("L138").Value = 7756.950000000001
$ws.Range("M138").Value = 942.8968000000004
$ws.Range("N138").Value = -18036.95

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1349
$ws.Range("I2").Value = 827.5714
$ws.Range("K2").Value = 827.5714
$ws.Range("M2").Value = -714.5714
$ws.Range("H32").Value = 3049.1143
$ws.Range("I32").Value = 3493.5356
$ws.Range("J32").Value = 1271.4286
$ws.Range("K32").Value = 3493.5356
$ws.Range("L32").Value = 1271.4286
$ws.Range("M32").Value = -3206.5356
$ws.Range("N32").Value = -1845.4286
$ws.Range("H45").Value = 1667.2941
$ws.Range("I45").Value = 1592.5
$ws.Range("J45").Value = 2016.3334
$ws.Range("K45").Value = 1592.5
$ws.Range("L45").Value = 2016.3334
$ws.Range("M45").Value = -1215.5
$ws.Range("N45").Value = -2770.3334
$ws.Range("H61").Value = 2323
$ws.Range("I61").Value = 2323
$ws.Range("K61").Value = 2323
$ws.Range("M61").Value = -2111
$ws.Range("H63").Value = 133353550
$ws.Range("I63").Value = 500020000
$ws.Range("K63").Value = 500020000
$ws.Range("M63").Value = -500019314
$ws.Range("H66").Value = 133353550
$ws.Range("I66").Value = 500020000
$ws.Range("K66").Value = 2500100000
$ws.Range("M66").Value = -2500096568
$ws.Range("H74").Value = 4314.231
$ws.Range("I74").Value = 4477.1304
$ws.Range("K74").Value = 4477.1304
$ws.Range("M74").Value = -3603.1304
$ws.Range("H77").Value = 4314.231
$ws.Range("I77").Value = 4477.1304
$ws.Range("K77").Value = 22385.652
$ws.Range("M77").Value = -18017.652
$ws.Range("H88").Value = 8334263
$ws.Range("I88").Value = 15152240
$ws.Range("K88").Value = 15152240
$ws.Range("M88").Value = -15151834
$ws.Range("H91").Value = 8334263
$ws.Range("I91").Value = 15152240
$ws.Range("K91").Value = 15152240
$ws.Range("M91").Value = -15150836
$ws.Range("H116").Value = 1349
$ws.Range("I116").Value = 827.5714
$ws.Range("K116").Value = 827.5714
$ws.Range("M116").Value = 1466.4286
$ws.Range("H119").Value = 28500
$ws.Range("J119").Value = 28500
$ws.Range("L119").Value = 28500
$ws.Range("N119").Value = -38176
$ws.Range("H122").Value = 3221.087
$ws.Range("I122").Value = 2705.6667
$ws.Range("J122").Value = 4187.5
$ws.Range("K122").Value = 8117.000100000001
$ws.Range("L122").Value = 12562.5
$ws.Range("M122").Value = -5667.000100000001
$ws.Range("N122").Value = -17462.5
$ws.Range("H125").Value = 42276.625
$ws.Range("J125").Value = 42276.625
$ws.Range("L125").Value = 42276.625
$ws.Range("N125").Value = -52116.625
$ws.Range("H136").Value = 2323
$ws.Range("I136").Value = 2323
$ws.Range("K136").Value = 6969
$ws.Range("M136").Value = -4419

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1349
$ws.Range("I3").Value = 827.5714
$ws.Range("K3").Value = 827.5714
$ws.Range("M3").Value = -713.5714
$ws.Range("H134").Value = 1401.5
$ws.Range("I134").Value = 1282
$ws.Range("J134").Value = 1999
$ws.Range("K134").Value = 3846
$ws.Range("L134").Value = 5997
$ws.Range("M134").Value = -1311
$ws.Range("N134").Value = -11067

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H31").Value = 5770.378
$ws.Range("I31").Value = 21184
$ws.Range("J31").Value = 3399.0513
$ws.Range("K31").Value = 21184
$ws.Range("L31").Value = 3399.0513
$ws.Range("M31").Value = -20889
$ws.Range("N31").Value = -3989.0513
$ws.Range("H34").Value = 5770.378
$ws.Range("I34").Value = 21184
$ws.Range("J34").Value = 3399.0513
$ws.Range("K34").Value = 21184
$ws.Range("L34").Value = 3399.0513
$ws.Range("M34").Value = -20982
$ws.Range("N34").Value = -3803.0513
$ws.Range("H105").Value = 3887.1428
$ws.Range("I105").Value = 3466.3333
$ws.Range("J105").Value = 4202.75
$ws.Range("K105").Value = 3466.3333
$ws.Range("L105").Value = 4202.75
$ws.Range("M105").Value = -1719.3333
$ws.Range("N105").Value = -7696.75
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 15378449
$ws.Range("I4").Value = 17574656
$ws.Range("J4").Value = 5002
$ws.Range("K4").Value = 52723968
$ws.Range("L4").Value = 15006
$ws.Range("M4").Value = -52723856
$ws.Range("N4").Value = -15230
$ws.Range("H12").Value = 334.1
$ws.Range("J12").Value = 335.16666
$ws.Range("L12").Value = 1005.49998
$ws.Range("N12").Value = -1351.49998
$ws.Range("H29").Value = 78.111115
$ws.Range("I29").Value = 88.2
$ws.Range("K29").Value = 264.6
$ws.Range("M29").Value = 12.39999999999998
$ws.Range("H56").Value = 19998.5
$ws.Range("I56").Value = 19998.5
$ws.Range("K56").Value = 19998.5
$ws.Range("M56").Value = -19468.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 5000
$ws.Range("J23").Value = 5000
$ws.Range("L23").Value = 5000
$ws.Range("N23").Value = -5446
$ws.Range("H122").Value = 2374.7273
$ws.Range("I122").Value = 2121.5334
$ws.Range("K122").Value = 6364.600199999999
$ws.Range("M122").Value = -3914.600199999999
$ws.Range("H126").Value = 8618.706
$ws.Range("J126").Value = 3499.8333
$ws.Range("L126").Value = 10499.4999
$ws.Range("N126").Value = -15439.4999
$ws.Range("H132").Value = 5137.654
$ws.Range("I132").Value = 3567
$ws.Range("J132").Value = 8104.4443
$ws.Range("K132").Value = 10701
$ws.Range("L132").Value = 24313.3329
$ws.Range("M132").Value = -8171
$ws.Range("N132").Value = -29373.3329

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3481.0625
$ws.Range("I7").Value = 2899.7273
$ws.Range("K7").Value = 2899.7273
$ws.Range("M7").Value = -2787.7273
$ws.Range("H40").Value = 3088.9565
$ws.Range("I40").Value = 2602.7646
$ws.Range("K40").Value = 2602.7646
$ws.Range("M40").Value = -2466.7646
$ws.Range("H46").Value = 2201.7036
$ws.Range("I46").Value = 1674.8235
$ws.Range("K46").Value = 1674.8235
$ws.Range("M46").Value = -1486.8235
$ws.Range("H55").Value = 420.17856
$ws.Range("I55").Value = 300.6842
$ws.Range("J55").Value = 672.44446
$ws.Range("K55").Value = 300.6842
$ws.Range("L55").Value = 672.44446
$ws.Range("M55").Value = -127.6842
$ws.Range("N55").Value = -1018.44446
$ws.Range("H122").Value = 5405.316
$ws.Range("I122").Value = 3700.182
$ws.Range("K122").Value = 11100.546
$ws.Range("M122").Value = -8650.545999999998
$ws.Range("H126").Value = 3481.0625
$ws.Range("I126").Value = 2899.7273
$ws.Range("K126").Value = 8699.1819
$ws.Range("M126").Value = -6229.1819
$ws.Range("H136").Value = 2498.3333
$ws.Range("I136").Value = 1986.4546
$ws.Range("K136").Value = 5959.3638
$ws.Range("M136").Value = -3409.3638

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H114").Value = 64999
$ws.Range("J114").Value = 64999
$ws.Range("L114").Value = 64999
$ws.Range("N114").Value = -73677
$ws.Range("H122").Value = 2186.5715
$ws.Range("I122").Value = 2087.2856
$ws.Range("K122").Value = 6261.8568
$ws.Range("M122").Value = -3811.8568

Write-Output "Updated cells: 224 set, 2 cleared."
